# Automatische test-sync: 2025-08-02 00:10:50
# Adds the new "Testmail #16" log row to the Logs sheet, extends the
# conditional-formatting ranges to cover it, and updates the Dashboard
# category counts (Intern verzoek / Actie voor medewerker now has 2 hits,
# ahead of Retour / Terugbetaling with 1).

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new row (row 6) to the Logs sheet -----------------------
$newRow = @(
    "Wil je dit even doorsturen?",
    "mailmind.test@zohomail.eu",
    "Testmail #16: Wil je dit even doorsturen?",
    "Intern verzoek / Actie voor medewerker",
    "Beste afzender,`nBedankt voor je bericht. Het is helaas niet helemaal duidelijk wat je precies bedoelt met `"Testmail #16: Wil je dit even doorsturen?`". Zou je meer informatie kunnen geven over wat je wilt dat er wordt doorgestuurd en naar welk e-mailadres dit moet gebeuren?`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent",
    "2025-08-02 00:10:07",
    "Ja",
    "Nee",
    "Ja",
    "Nee"
)

for ($col = 1; $col -le $newRow.Length; $col++) {
    $logs.Cells.Item(6, $col).Value = $newRow[$col - 1]
}

# Writing the long, multi-line reply text can trigger an implicit
# "custom row height" on entry; restore the default (auto) row height so
# row 6 matches the rest of the sheet (no explicit ht/customHeight).
$logs.Rows.Item(6).AutoFit()

# --- 2. Extend conditional formatting ranges from row 5 to row 6 -----------
$logs.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))
$logs.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))
$logs.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H6"))
$logs.Range("I2:I6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I6"))
$logs.Range("J2:J6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J6"))

# --- 3. Update the dashboard summary (row order swap + new count) ----------
$dash.Cells.Item(3, 1).Value = "Intern verzoek / Actie voor medewerker"
$dash.Cells.Item(3, 2).Value = 2
$dash.Cells.Item(4, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(4, 2).Value = 1
